$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new rows before row 7 (i.e. at rows 5-6), shifting the old
# rows 5-9 (Manage Wait times.., Test Execution Path.., Extent Report Config..,
# version type.., Application URL-DUP..) down to rows 7-11. This preserves
# their original cell styles/content automatically. ---
$ws.Range("A5:A6").EntireRow.Insert()

# --- Row 4: "Extent Report Configuration" label moves away (now blank,
# format kept) since it no longer maps 1:1 with "Execution Type" ---
$ws.Range("A4").ClearContents()

# --- Rows 5 and 6 are new blank rows: restore the left-column border
# styling (copy format from A4, which already carries that style) and
# fill in the new "ParallelTests" / "Execution env" config rows ---
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A4").Copy()
$ws.Range("A6").PasteSpecial(-4122)

$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B4").Copy()
$ws.Range("B6").PasteSpecial(-4122)

$ws.Range("C4").Copy()
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("C4").Copy()
$ws.Range("C6").PasteSpecial(-4122)

$ws.Range("B5").Value = "ParallelTests"
$ws.Range("C5").Value = 3

$ws.Range("B6").Value = "Execution env"
$ws.Range("C6").Value = "Local"

# --- New data validation rules for the two new cells ---
$ws.Range("C5").Validation.Add(3, 1, 1, '"1,2,3,4,5,6,7,8"')
$ws.Range("C6").Validation.Add(3, 1, 1, '"Local, Remote"')

# --- The external hyperlink that used to sit on the "TestExeFile" row
# (old C6) now lives on the same row after it shifted down to C8 ---
$target = $ws.Range("C8").Value()
$ws.Range("C6").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C8"), "file:///\\data\TestExecution.xlsx")

# --- Grow Table1 to cover the two extra rows ---
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:C11"))

# --- Restore the active selection that was left in the source file ---
$ws.Range("H5").Select()
